{"js": "// The document contains a single table laid out as repeating groups of\n// 4 rows (1 data row with 5 math-fact cells, followed by 3 blank rows).\n// This script rewrites the 25 populated cells (5 data rows x 5 columns)\n// with their new values, addressed positionally by (row, column) so that\n// duplicate / reused text values (e.g. \"76\u00f79=8, 4\" appears twice, and\n// \"59\u00f72=29, 1\" is both a target and, later, a source) are never confused\n// with each other the way a global find-and-replace would.\nconst edits = [\n  { row: 0, col: 0, oldText: \"27\u00f75=5, 2\", newText: \"55\u00f73=18, 1\" },\n  { row: 0, col: 1, oldText: \"87\u00f77=12, 3\", newText: \"30\u00f76=5, 0\" },\n  { row: 0, col: 2, oldText: \"21\u00f72=10, 1\", newText: \"82\u00f79=9, 1\" },\n  { row: 0, col: 3, oldText: \"76\u00f79=8, 4\", newText: \"13\u00f76=2, 1\" },\n  { row: 0, col: 4, oldText: \"11\u00f73=3, 2\", newText: \"50\u00f76=8, 2\" },\n  { row: 4, col: 0, oldText: \"92\u00f77=13, 1\", newText: \"73\u00f78=9, 1\" },\n  { row: 4, col: 1, oldText: \"76\u00f79=8, 4\", newText: \"49\u00f73=16, 1\" },\n  { row: 4, col: 2, oldText: \"37\u00f73=12, 1\", newText: \"81\u00f72=40, 1\" },\n  { row: 4, col: 3, oldText: \"18\u00f75=3, 3\", newText: \"55\u00f72=27, 1\" },\n  { row: 4, col: 4, oldText: \"70\u00f79=7, 7\", newText: \"59\u00f72=29, 1\" },\n  { row: 8, col: 0, oldText: \"48\u00f76=8, 0\", newText: \"63\u00f77=9, 0\" },\n  { row: 8, col: 1, oldText: \"59\u00f72=29, 1\", newText: \"11\u00f78=1, 3\" },\n  { row: 8, col: 2, oldText: \"98\u00f79=10, 8\", newText: \"86\u00f79=9, 5\" },\n  { row: 8, col: 3, oldText: \"86\u00f73=28, 2\", newText: \"28\u00f79=3, 1\" },\n  { row: 8, col: 4, oldText: \"48\u00f79=5, 3\", newText: \"13\u00f77=1, 6\" },\n  { row: 12, col: 0, oldText: \"37\u00f76=6, 1\", newText: \"37\u00f77=5, 2\" },\n  { row: 12, col: 1, oldText: \"21\u00f78=2, 5\", newText: \"78\u00f73=26, 0\" },\n  { row: 12, col: 2, oldText: \"51\u00f72=25, 1\", newText: \"69\u00f76=11, 3\" },\n  { row: 12, col: 3, oldText: \"16\u00f73=5, 1\", newText: \"36\u00f74=9, 0\" },\n  { row: 12, col: 4, oldText: \"36\u00f72=18, 0\", newText: \"14\u00f77=2, 0\" },\n  { row: 16, col: 0, oldText: \"90\u00f78=11, 2\", newText: \"86\u00f77=12, 2\" },\n  { row: 16, col: 1, oldText: \"56\u00f73=18, 2\", newText: \"46\u00f77=6, 4\" },\n  { row: 16, col: 2, oldText: \"99\u00f77=14, 1\", newText: \"63\u00f72=31, 1\" },\n  { row: 16, col: 3, oldText: \"90\u00f74=22, 2\", newText: \"87\u00f73=29, 0\" },\n  { row: 16, col: 4, oldText: \"95\u00f74=23, 3\", newText: \"10\u00f79=1, 1\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Load every target cell's current text so we can sanity-check it against\n// the expected \"before\" value prior to overwriting it.\nconst cells = edits.map((e) => rows.items[e.row].cells.items[e.col]);\ncells.forEach((cell) => cell.body.load(\"text\"));\nawait context.sync();\n\nfor (let i = 0; i < edits.length; i++) {\n  const cell = cells[i];\n  const expected = edits[i].oldText;\n  const actual = (cell.body.text || \"\").trim();\n  if (actual !== expected) {\n    throw new Error(\n      `Cell (row ${edits[i].row}, col ${edits[i].col}) text mismatch: ` +\n        `expected \"${expected}\" but found \"${actual}\"`\n    );\n  }\n  cell.value = edits[i].newText;\n}\n\nawait context.sync();\n", "ps1": "# The document contains a single table laid out as repeating groups of\n# 4 rows (1 data row with 5 math-fact cells, followed by 3 blank rows).\n# This script rewrites the 25 populated cells (5 data rows x 5 columns)\n# with their new values, addressed positionally by (row, column) - using\n# Word's 1-based Table.Cell(row, col) addressing - so that duplicate /\n# reused text values (e.g. \"76\u00f79=8, 4\" appears twice, and \"59\u00f72=29, 1\"\n# is both a target and, later, a source) are never confused with each\n# other the way a global Find/Replace would be.\n\n$edits = @(\n    @{ Row = 1; Col = 1; OldText = \"27\u00f75=5, 2\"; NewText = \"55\u00f73=18, 1\" },\n    @{ Row = 1; Col = 2; OldText = \"87\u00f77=12, 3\"; NewText = \"30\u00f76=5, 0\" },\n    @{ Row = 1; Col = 3; OldText = \"21\u00f72=10, 1\"; NewText = \"82\u00f79=9, 1\" },\n    @{ Row = 1; Col = 4; OldText = \"76\u00f79=8, 4\"; NewText = \"13\u00f76=2, 1\" },\n    @{ Row = 1; Col = 5; OldText = \"11\u00f73=3, 2\"; NewText = \"50\u00f76=8, 2\" },\n    @{ Row = 5; Col = 1; OldText = \"92\u00f77=13, 1\"; NewText = \"73\u00f78=9, 1\" },\n    @{ Row = 5; Col = 2; OldText = \"76\u00f79=8, 4\"; NewText = \"49\u00f73=16, 1\" },\n    @{ Row = 5; Col = 3; OldText = \"37\u00f73=12, 1\"; NewText = \"81\u00f72=40, 1\" },\n    @{ Row = 5; Col = 4; OldText = \"18\u00f75=3, 3\"; NewText = \"55\u00f72=27, 1\" },\n    @{ Row = 5; Col = 5; OldText = \"70\u00f79=7, 7\"; NewText = \"59\u00f72=29, 1\" },\n    @{ Row = 9; Col = 1; OldText = \"48\u00f76=8, 0\"; NewText = \"63\u00f77=9, 0\" },\n    @{ Row = 9; Col = 2; OldText = \"59\u00f72=29, 1\"; NewText = \"11\u00f78=1, 3\" },\n    @{ Row = 9; Col = 3; OldText = \"98\u00f79=10, 8\"; NewText = \"86\u00f79=9, 5\" },\n    @{ Row = 9; Col = 4; OldText = \"86\u00f73=28, 2\"; NewText = \"28\u00f79=3, 1\" },\n    @{ Row = 9; Col = 5; OldText = \"48\u00f79=5, 3\"; NewText = \"13\u00f77=1, 6\" },\n    @{ Row = 13; Col = 1; OldText = \"37\u00f76=6, 1\"; NewText = \"37\u00f77=5, 2\" },\n    @{ Row = 13; Col = 2; OldText = \"21\u00f78=2, 5\"; NewText = \"78\u00f73=26, 0\" },\n    @{ Row = 13; Col = 3; OldText = \"51\u00f72=25, 1\"; NewText = \"69\u00f76=11, 3\" },\n    @{ Row = 13; Col = 4; OldText = \"16\u00f73=5, 1\"; NewText = \"36\u00f74=9, 0\" },\n    @{ Row = 13; Col = 5; OldText = \"36\u00f72=18, 0\"; NewText = \"14\u00f77=2, 0\" },\n    @{ Row = 17; Col = 1; OldText = \"90\u00f78=11, 2\"; NewText = \"86\u00f77=12, 2\" },\n    @{ Row = 17; Col = 2; OldText = \"56\u00f73=18, 2\"; NewText = \"46\u00f77=6, 4\" },\n    @{ Row = 17; Col = 3; OldText = \"99\u00f77=14, 1\"; NewText = \"63\u00f72=31, 1\" },\n    @{ Row = 17; Col = 4; OldText = \"90\u00f74=22, 2\"; NewText = \"87\u00f73=29, 0\" },\n    @{ Row = 17; Col = 5; OldText = \"95\u00f74=23, 3\"; NewText = \"10\u00f79=1, 1\" }\n)\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\nforeach ($edit in $edits) {\n    $cell = $table.Cell($edit.Row, $edit.Col)\n    $range = $cell.Range\n    # Cell.Range.Text carries the trailing end-of-cell marker (CR + BEL);\n    # strip it before comparing against the expected \"before\" value.\n    $actual = $range.Text.TrimEnd([char]13, [char]7)\n    if ($actual -ne $edit.OldText) {\n        throw \"Cell (row $($edit.Row), col $($edit.Col)) text mismatch: expected '$($edit.OldText)' but found '$actual'\"\n    }\n    $range.Text = $edit.NewText\n}\n"}
